$d = $word.ActiveDocument

# 1. Insert "(aka kmeans, k means)" after "k-means clustering" and before
#    " is a method of vector quantization..."
$d.Content.Find.Execute(
    "k-means clustering is a method",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "k-means clustering (aka kmeans, k means) is a method",
    2)

# 2. Simplify "Voronoi cells." back into plain, unmarked text.
$d.Content.Find.Execute(
    "Voronoi cells.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Voronoi cells.",
    2)

# 3. Simplify "Rocchio algorithm." back into plain, unmarked text.
$d.Content.Find.Execute(
    "Rocchio algorithm.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Rocchio algorithm.",
    2)
